# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
# Numeric-looking Price (column D) values must stay TEXT (the source stores them
# as inline strings), so they are written with a leading apostrophe to force
# Excel to keep them as text instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.209.79'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '1.654.62'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('D5').Value = '''218.47'
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('D6').Value = '''0.5200'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('D9').Value = '''0.06311'
$ws.Range('E9').Value = '  -1.27%  '
$ws.Range('D10').Value = '''21.27'
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('D11').Value = '''0.07737'
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').Value = '1.649.81'
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').Value = '''4.424'
$ws.Range('D14').Value = '''0.5457'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').Value = '0.0₅8197'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '''64.73'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '26.197.55'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('D19').Value = '''4.677'
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('D20').Value = '''191.07'
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').Value = '''10.16'
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('D22').Value = '''6.184'
$ws.Range('E22').Value = '  -2.96%  '
$ws.Range('D23').Value = '''1.006'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('D24').Value = '''138.58'
$ws.Range('E24').Value = '  -3.17%  '
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').Value = '''7.285'
$ws.Range('E26').Value = '  -1.77%  '
$ws.Range('D27').Value = '''16.06'
$ws.Range('E27').Value = '  +0.38%  '
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('D29').Value = '''0.06061'
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('D30').Value = '''1.284'
$ws.Range('E30').Value = '  +1.42%  '
$ws.Range('D31').Value = '''3.547'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').Value = '''3.359'
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('D33').Value = '''1.652'
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('D34').Value = '''0.9848'
$ws.Range('E34').Value = '  -1.50%  '
$ws.Range('D35').Value = '''2.412'
$ws.Range('E35').Value = '  +0.46%  '
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('D37').Value = '''0.5925'
$ws.Range('E37').Value = '  +4.61%  '
$ws.Range('D38').Value = '''0.01596'
$ws.Range('E38').Value = '  -0.63%  '
$ws.Range('D39').Value = '''5.952'
$ws.Range('E39').Value = '  +0.69%  '
$ws.Range('D40').Value = '''0.8628'
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('D41').Value = '1.057.21'
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '1.796.28'
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '''57.37'
$ws.Range('E45').Value = '  +2.23%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('D47').Value = '''1.004'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '''8.041'
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('D49').Value = '''0.05177'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').Value = '''1.466'
$ws.Range('E50').Value = '  +4.94%  '
$ws.Range('D51').Value = '''0.4232'
$ws.Range('E51').Value = '  +0.42%  '
